$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): columns reordered / relabeled ---
# Only the cells whose text actually differs from what is already there need touching.
$ws.Range("A1").Value2 = "Nhóm KH"
$ws.Range("D1").Value2 = "DVT - KH"
$ws.Range("E1").Value2 = "Số lượng - KH"

# --- Data row (row 2): new single remaining data row ---
$ws.Range("A2").Value2 = "Emart"
$ws.Range("B2").Value2 = 8936040451101
$ws.Range("C2").Value2 = 8936040451101
$ws.Range("D2").Value2 = "CAI"
$ws.Range("E2").Value2 = 1

# D2 and E2 lose the thin-box border that the old quantity/unit cells had
$ws.Range("D2:E2").Borders.LineStyle = 0

# --- Remove the now-unused extra data rows (old rows 3 and 4) ---
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()

# --- Column widths for B and E ---
$ws.Columns.Item(2).ColumnWidth = 12.17
$ws.Columns.Item(5).ColumnWidth = 12

# --- Selection moves to G6 ---
$ws.Range("G6").Select()
